# Weekly price update: insert a new record as row 61, pushing the existing
# rows 61-165 down to 62-166 (dimension grows from A1:R165 to A1:R166).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 61, shifting everything below it down by one.
$ws.Rows("61:61").Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A61").Value = 10
$ws.Range("B61").Value = "Vega Modelo de Temuco"
$ws.Range("C61").Value = "La Araucanía"
$ws.Range("D61").Value = 44771
$ws.Range("E61").Value = 9
$ws.Range("F61").Value = 100112012
$ws.Range("G61").Value = "Espinaca"
$ws.Range("H61").Value = "Sin especificar"
$ws.Range("I61").Value = "Primera"
$ws.Range("J61").Value = 20
$ws.Range("K61").Value = 14000
$ws.Range("L61").Value = 14000
$ws.Range("M61").Value = 14000
$ws.Range("N61").Value = "$/docena de atados"
$ws.Range("O61").Value = "Región de La Araucanía"
$ws.Range("P61").Value = 4667
$ws.Range("Q61").Value = 3
$ws.Range("R61").Value = "Hortaliza"
